# Resume update script
# Applies the edits described in the commit "Fix 4 insights and database management"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace text once within the whole document (text is assumed unique)
# ---------------------------------------------------------------------------
function Replace-Unique($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# Helper: replace text once, scoped to a specific paragraph (by 1-based index)
# ---------------------------------------------------------------------------
function Replace-InParagraph($index, $old, $new) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ===========================================================================
# 1. TEXT-ONLY REPLACEMENTS (done first, while paragraph indices still match
#    the original document layout)
# ===========================================================================

# --- Summary paragraph -----------------------------------------------------
Replace-Unique `
    "Soy una persona proactiva y emprendedora, con gran facilidad de aprendizaje. Alta capacidad de trabajo en equipo y de amoldarse ante cualquier situación y poder generar soluciones ante cualquier imprevisto. Tengo iniciativa propia para generar ideas que ayuden a la rapidez del trabajo y muchas ganas de superación." `
    "Ingeniero de Sistemas con +4 años de experiencia en desarrollado de software. Me apasiona la tecnología, inicié como Freelance para proyectos Web en Java, continué como Programador Java Junior en DIRESA, Programador GeneXus Java y C++ en una Consultora de Software, he desarrollado soluciones contables para bancos en Colombia, además he desarrollado soluciones web de facturación electrónica para empresas privadas y para la plataforma web de SUNAT. Tengo sólidos conocimientos en Node.js, Java, Golang, Vue.js, ReactJs, Angular, TypeScript, SaaS, Cloud. Actualmente, soy Arquitecto de Software en Interseguro, donde continúo ampliando mis habilidades técnicas y de gestión."

# --- FORMACIÓN ---------------------------------------------------------------
Replace-Unique `
    "Ingeniero Informático at Universidad Ricardo Palma, 2016-2020, Lima, Peru" `
    "Ingeniería de Sistemas, Universidad Nacional José Faustino Sánchez Carrión, 2015-2020"

# --- CURSOS Y CERTIFICACIONES ----------------------------------------------
Replace-Unique "Angular Avanzando UDEMY" "Programador Java - SAP Integral, SAP HANNA, SAP ERP, Universidad Nacional de Ingeniería"
Replace-Unique "Diseño Web Avanzado UDEMY" "Especialización Java Web Developer - Spring Boot 2, NodeJs, Galaxy Training"
# (paragraphs for "AWS Cloud Practitioner UDEMY", "REACT Avanzando UDEMY" and
#  "Diseño web Avanzando UDEMY" are removed entirely further below)

# --- TECNOLOGÍAS -------------------------------------------------------------
Replace-Unique "Angular 5+" "JS/Angular/VueJs/React"
Replace-Unique "Jenkins" "Java/C#/C++/Nodejs/Golang"
Replace-Unique "Bitbucket" "GoogleCloudPlatform/Azure"
Replace-Unique "Jira" "MongoDB/MySQL/DB2/SQLS"
# "AWS" is ambiguous (also appears inside "AWS Cloud Practitioner UDEMY"),
# so scope the replace to its specific paragraph.
Replace-InParagraph 16 "AWS" "Postman/SoapUI"
Replace-Unique "Jasmine" "SonarQube/Testing"
Replace-Unique "Karma" "Maven/Graddle"
Replace-Unique ".Net Core" "Spring Framework, S. Data, S. Security, S.Cloud"
Replace-Unique "Spring JDBC" "Firebase9, Identity Platform"
# (paragraphs "React.js", "Boostrap", "Scrum", "Azure", "SQL Server", "Oracle",
#  "Microfrontend", "Java", "Springboot" are removed entirely further below)

# --- EXPERIENCIA -------------------------------------------------------------

# Job 1 (paragraphs 31-33)
Replace-InParagraph 31 "Desarrollador FullStack" "Arquitecto de Software"
Replace-InParagraph 32 "NTT DATA S.A.C." "Interseguro S.A"
Replace-InParagraph 32 "DEC 2020 - ACTUALIDAD" "Apr 2023 - Present"
Replace-InParagraph 33 `
    "Desarrollo de micro frontends y microservicios en java orientados al flujo bancario del BCP." `
    "Responsable del diseño y desarrollo de aplicaciones, empleando diversas tecnologías, como Node.js, Golang y Java, así como frameworks frontend como Vue (Nuxt3), React (NextJs) y Angular. Uso metodologías como Integración Continua, TDD y Desarrollo Ágil. En el backend, sigo una arquitectura de microservicios, incluyendo microservicios orquestadores. Para el frontend, adopto arquitecturas hexagonales y de servicios, además de diseñar y crear microfrontends en Vue 3 con la biblioteca single-spa-vue para integraciones en aplicaciones de una sola página (SPA). Entre mis responsabilidades se encuentra la creación de arquetipos y directrices, abarcando desde la definición de diagramas de arquitectura hasta la elaboración de estándares de desarrollo y buenas prácticas. Defino estructuras de directorios, gestiono integraciones con bases de datos (MongoDB, PostgreSQL, Oracle) y uso ORMs e integraciones REST API. Implemento registro de eventos personalizados (loggers) y reportes de salud (liveness, readiness), implementación Google Auth con Identity Platform en el front, generación de imágenes Docker, integración continua con GitLab Pipelines, despliegue en entornos como Cloud Run y Kubernetes, y documentación con Swagger."

# Job 2 (paragraphs 34-36)
Replace-InParagraph 34 "Desarrollador FULLSTACK" "Programador Sr Java"
Replace-InParagraph 35 "SIGCOMT S.A.C." "TCI S.A"
Replace-InParagraph 35 "JUL 2020 - DEC 2020" "Jan 2022 - Mar 2023"
Replace-InParagraph 36 `
    "Programación de servicios Rest en .Net y flujos frontend en Angular para la entidad Calida." `
    "Responsable del análisis y desarrollo de funcionalidades para las aplicaciones web de facturación electrónica y OSE, Suite Movil y SuiteOnline, ambas aplicaciones de control de Guías de Remisión, de Facturación Electrónica en los sistemas de gestión ePortal y eGestor, transmisión de Comprobantes con manejo de Colas con RabbitMQ.Entre los más, Gestión de Usuarios, Dispositivos, Empresas, Guías, Comprobantes, Inventarios, Tipo de cambio, Reportes, etc."

# Job 3 (paragraphs 37-39)
Replace-InParagraph 37 "Desarrollador FULLSTACK" "Analista Programador Sr"
Replace-InParagraph 38 "STAMINER S.A.C." "Canvia"
Replace-InParagraph 38 "JUL 2019 - JUL 2020" "Agos 2022 - Feb 2023"
Replace-InParagraph 39 `
    "Programación de servicios Rest y web app en angular para una red social de deportistas." `
    "Desarrollo de mejoras y nuevas funciones para la plataforma web SUNAT y otras empresas privadas bajo la práctica de Integración Continua, participé en el equipo de desarrollo de Microservicios, para el desarrollo de la funcionalidad de Registros de Diligencia Inicial y Diligencia de Tránsito para Declaraciones de Transito Aduanero Internacional SUNAT y consultas de comprobantes SUNAT."

# Job 4 (paragraphs 40-42)
Replace-InParagraph 40 "Desarrollador FULLSTACK" "Programador Semi Sr"
Replace-InParagraph 41 "GESTION Y SISTEMAS S.A.C." "Hiper SA"
Replace-InParagraph 41 "JUL 2018 - JUL 2019" "Jan 2022 - Jul 2022"
Replace-InParagraph 42 `
    "Programación de servicios Rest para su implementación, también desarrollo en framework Angular para proyectos de los clientes Exalmar Pesquera, Ripley Bodegas y La Contraloria." `
    "Desarrollo y mantenimiento de procesos y transacciones para bancos y servicios electrónicos en Colombia, Depuración de vulnerabilidades, modificación de fuentes, mantenimiento y ejecución de disparadores, manejo de repositorios, mantenimiento retiro consultas OTP y Pruebas unitarias."

# ===========================================================================
# 2. DELETIONS (remove obsolete bullet paragraphs). Higher-index block first
#    so the as-yet-unprocessed lower-index block keeps its original numbering.
# ===========================================================================

# Remove technology bullets: React.js, Boostrap, Scrum, Azure, SQL Server,
# Oracle, Microfrontend, Java, Springboot (paragraphs 20-28)
$startP = $d.Paragraphs.Item(20)
$endP = $d.Paragraphs.Item(28)
$d.Range($startP.Range.Start, $endP.Range.End).Delete() | Out-Null

# Remove course bullets: AWS Cloud Practitioner UDEMY, REACT Avanzando UDEMY,
# Diseño web Avanzando UDEMY (paragraphs 8-10)
$startP2 = $d.Paragraphs.Item(8)
$endP2 = $d.Paragraphs.Item(10)
$d.Range($startP2.Range.Start, $endP2.Range.End).Delete() | Out-Null

# ===========================================================================
# 3. INSERTIONS: append the two new job entries at the end of the body
# ===========================================================================

$endRange = $d.Content
$endRange.Collapse(0) | Out-Null

$newJobsXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships>
</pkg:xmlData></pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Programador Genexus Jr</w:t></w:r></w:p>
<w:p><w:r><w:t>PROEMSA SAC</w:t><w:br/><w:t>Sep 2020 - Dec 2021</w:t></w:r></w:p>
<w:p><w:r><w:t>Desarrollo de aplicaciones web a medida para empresas como (INDUSFER, INCALPACA, LALENNA, SAN GABAN, etc).</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Programador Java Jr</w:t></w:r></w:p>
<w:p><w:r><w:t>DIRESA Lima</w:t><w:br/><w:t>Dec 2019 - Mar 2020</w:t></w:r></w:p>
<w:p><w:r><w:t>Diseño web de la página web DIRESA, mantenimiento y desarrollo de Reportes documentación de Serums y mantenimiento del Portal de Transparencia.</w:t></w:r></w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

$endRange.InsertXML($newJobsXml) | Out-Null

# ===========================================================================
# 4. HEADER
# ===========================================================================

$hdr = $d.Sections.Item(1).Headers.Item(1)
$hdr.Range.Find.Execute("Harold Portillo", $true, $true, $false, $false, $false, $true, 1, $false, "Angelo Lugo", 2) | Out-Null
$hdr.Range.Find.Execute("Lima, Peru", $true, $true, $false, $false, $false, $true, 1, $false, "Lima, Perú", 2) | Out-Null
$hdr.Range.Find.Execute("Peru", $true, $true, $false, $false, $false, $true, 1, $false, "Perú", 2) | Out-Null
